$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Prix Spot": insert a new day column ("18-dec") before column EU,
# shifting every column from EU onward one place to the right
# (EU -> EV ... FY -> FZ). Dimension grows from A1:FY25 to A1:FZ25.
# ----------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Insert a new blank column at EU (existing EU:FY shift to FV:FZ)
$wsSpot.Range("EU1:EU25").EntireColumn.Insert()

# New header cell for the inserted day
$wsSpot.Range("EU1").Value = "18-dec"

# New data cells for the inserted day (no data available yet -> "-")
$wsSpot.Range("EU2:EU25").Value = "-"

# ----------------------------------------------------------------------
# Sheet "Gaz": append a new row with the next day's data
# ----------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A181").NumberFormat = "@"
$wsGaz.Range("A181").Value = "2025-12-16"
$wsGaz.Range("A181").ClearFormats()
$wsGaz.Range("B181").Value = 25.475

# ----------------------------------------------------------------------
# Sheet "CO2": append a new row with the next day's data
# ----------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A181").NumberFormat = "@"
$wsCo2.Range("A181").Value = "2025-12-16"
$wsCo2.Range("A181").ClearFormats()
$wsCo2.Range("B181").Value = 85.08
